# Rename the "sbml_models" sheet to "models"
$wb = $excel.ActiveWorkbook
$sbmlSheet = $wb.Worksheets.Item("sbml_models")
$sbmlSheet.Name = "models"

# Update reactor_intervals!H7 to reference the new json model instead of the inline formula text
$reactorSheet = $wb.Worksheets.Item("reactor_intervals")
$reactorSheet.Range("H7").Value = "Glucose_open_fermentation.json"

# Add a new row (row 7) to the "models" sheet describing the open_fermentation model
$modelsSheet = $wb.Worksheets.Item("models")

# Carry over the same fill/border formatting used by the rows above
# (columns B:D share one style, E:G share another).
$modelsSheet.Range("B6:D6").Copy()
$modelsSheet.Range("B7:D7").PasteSpecial(-4122)
$modelsSheet.Range("E6:G6").Copy()
$modelsSheet.Range("E7:G7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$modelsSheet.Range("F7").Value = "Propionate, Acetate, Butyrate"
$modelsSheet.Range("G7").Value = "prop, ace, but"
$modelsSheet.Range("C7").Value = "Glucose, pH"
$modelsSheet.Range("D7").Value = "glu,pH"
$modelsSheet.Range("A7").Value = "Glucose_open_fermentation.json"
$modelsSheet.Range("B7").Value = 0
$modelsSheet.Range("E7").Value = 0

# Match cell selections as left by the editing session.
# "reactor_intervals" stays the active sheet, so select there last.
$modelsSheet.Range("C15").Select()
$reactorSheet.Activate()
$reactorSheet.Range("G9").Select()
